# Apply NATMI LR-pair update (Col1a1-Itga2): recompute the 12 existing sender/target
# cluster combination rows and add the 4 missing "sCs" sender rows, per Dr Hou's advice.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a1"
$ws.Range("C2").Value = "Itga2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.32821866666667
$ws.Range("H2").Value = 60.984656
$ws.Range("I2").Value = 0.004181898474048532
$ws.Range("J2").Value = 0.004181898474048532
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.740822
$ws.Range("N2").Value = 5.222466
$ws.Range("O2").Value = 0.4863878955914668
$ws.Range("P2").Value = 0.4863878955914669
$ws.Range("Q2").Value = 35.387810275744
$ws.Range("R2").Value = 318.490292481696
$ws.Range("S2").Value = 0.002034024798369632
$ws.Range("T2").Value = 0.002034024798369632

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a1"
$ws.Range("C3").Value = "Itga2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.32821866666667
$ws.Range("H3").Value = 60.984656
$ws.Range("I3").Value = 0.004181898474048532
$ws.Range("J3").Value = 0.004181898474048532
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.077748
$ws.Range("N3").Value = 3.233244
$ws.Range("O3").Value = 0.3011241710513264
$ws.Range("P3").Value = 0.3011241710513265
$ws.Range("Q3").Value = 21.90869701156267
$ws.Range("R3").Value = 197.178273104064
$ws.Range("S3").Value = 0.001259270711418671
$ws.Range("T3").Value = 0.001259270711418671

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a1"
$ws.Range("C4").Value = "Itga2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.32821866666667
$ws.Range("H4").Value = 60.984656
$ws.Range("I4").Value = 0.004181898474048532
$ws.Range("J4").Value = 0.004181898474048532
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03488166666666666
$ws.Range("N4").Value = 0.104645
$ws.Range("O4").Value = 0.009745982326006345
$ws.Range("P4").Value = 0.009745982326006345
$ws.Range("Q4").Value = 0.7090821474577778
$ws.Range("R4").Value = 6.38173932712
$ws.Range("S4").Value = 0.0000407567086172299
$ws.Range("T4").Value = 0.0000407567086172299

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col1a1"
$ws.Range("C5").Value = "Itga2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.32821866666667
$ws.Range("H5").Value = 60.984656
$ws.Range("I5").Value = 0.004181898474048532
$ws.Range("J5").Value = 0.004181898474048532
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7256300000000001
$ws.Range("N5").Value = 2.17689
$ws.Range("O5").Value = 0.2027419510312003
$ws.Range("P5").Value = 0.2027419510312003
$ws.Range("Q5").Value = 14.75076531109334
$ws.Range("R5").Value = 132.75688779984
$ws.Range("S5").Value = 0.0008478462556429987
$ws.Range("T5").Value = 0.0008478462556429987

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a1"
$ws.Range("C6").Value = "Itga2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4809.896321333334
$ws.Range("H6").Value = 14429.688964
$ws.Range("I6").Value = 0.9894865072215304
$ws.Range("J6").Value = 0.9894865072215304
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.740822
$ws.Range("N6").Value = 5.222466
$ws.Range("O6").Value = 0.4863878955914668
$ws.Range("P6").Value = 0.4863878955914669
$ws.Range("Q6").Value = 8373.173333896137
$ws.Range("R6").Value = 75358.56000506523
$ws.Range("S6").Value = 0.4812742599636309
$ws.Range("T6").Value = 0.481274259963631

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a1"
$ws.Range("C7").Value = "Itga2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4809.896321333334
$ws.Range("H7").Value = 14429.688964
$ws.Range("I7").Value = 0.9894865072215304
$ws.Range("J7").Value = 0.9894865072215304
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.077748
$ws.Range("N7").Value = 3.233244
$ws.Range("O7").Value = 0.3011241710513264
$ws.Range("P7").Value = 0.3011241710513265
$ws.Range("Q7").Value = 5183.856140524357
$ws.Range("R7").Value = 46654.70526471922
$ws.Range("S7").Value = 0.2979583042535557
$ws.Range("T7").Value = 0.2979583042535557

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col1a1"
$ws.Range("C8").Value = "Itga2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4809.896321333334
$ws.Range("H8").Value = 14429.688964
$ws.Range("I8").Value = 0.9894865072215304
$ws.Range("J8").Value = 0.9894865072215304
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03488166666666666
$ws.Range("N8").Value = 0.104645
$ws.Range("O8").Value = 0.009745982326006345
$ws.Range("P8").Value = 0.009745982326006345
$ws.Range("Q8").Value = 167.7772001819756
$ws.Range("R8").Value = 1509.99480163778
$ws.Range("S8").Value = 0.009643518011202785
$ws.Range("T8").Value = 0.009643518011202785

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col1a1"
$ws.Range("C9").Value = "Itga2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4809.896321333334
$ws.Range("H9").Value = 14429.688964
$ws.Range("I9").Value = 0.9894865072215304
$ws.Range("J9").Value = 0.9894865072215304
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7256300000000001
$ws.Range("N9").Value = 2.17689
$ws.Range("O9").Value = 0.2027419510312003
$ws.Range("P9").Value = 0.2027419510312003
$ws.Range("Q9").Value = 3490.205067649108
$ws.Range("R9").Value = 31411.84560884197
$ws.Range("S9").Value = 0.2006104249931409
$ws.Range("T9").Value = 0.2006104249931409

# Row 10: M2 -> ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Col1a1"
$ws.Range("C10").Value = "Itga2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.69506
$ws.Range("H10").Value = 8.085180000000001
$ws.Range("I10").Value = 0.000554424737665286
$ws.Range("J10").Value = 0.000554424737665286
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.740822
$ws.Range("N10").Value = 5.222466
$ws.Range("O10").Value = 0.4863878955914668
$ws.Range("P10").Value = 0.4863878955914669
$ws.Range("Q10").Value = 4.69161973932
$ws.Range("R10").Value = 42.22457765388
$ws.Range("S10").Value = 0.0002696654814168695
$ws.Range("T10").Value = 0.0002696654814168696

# Row 11: M2 -> FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col1a1"
$ws.Range("C11").Value = "Itga2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.69506
$ws.Range("H11").Value = 8.085180000000001
$ws.Range("I11").Value = 0.000554424737665286
$ws.Range("J11").Value = 0.000554424737665286
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.077748
$ws.Range("N11").Value = 3.233244
$ws.Range("O11").Value = 0.3011241710513264
$ws.Range("P11").Value = 0.3011241710513265
$ws.Range("Q11").Value = 2.90459552488
$ws.Range("R11").Value = 26.14135972392
$ws.Range("S11").Value = 0.0001669506895398084
$ws.Range("T11").Value = 0.0001669506895398084

# Row 12: M2 -> M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col1a1"
$ws.Range("C12").Value = "Itga2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.69506
$ws.Range("H12").Value = 8.085180000000001
$ws.Range("I12").Value = 0.000554424737665286
$ws.Range("J12").Value = 0.000554424737665286
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.03488166666666666
$ws.Range("N12").Value = 0.104645
$ws.Range("O12").Value = 0.009745982326006345
$ws.Range("P12").Value = 0.009745982326006345
$ws.Range("Q12").Value = 0.09400818456666667
$ws.Range("R12").Value = 0.8460736611000002
$ws.Range("S12").Value = 0.000005403413694386582
$ws.Range("T12").Value = 0.000005403413694386582

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col1a1"
$ws.Range("C13").Value = "Itga2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.69506
$ws.Range("H13").Value = 8.085180000000001
$ws.Range("I13").Value = 0.000554424737665286
$ws.Range("J13").Value = 0.000554424737665286
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.7256300000000001
$ws.Range("N13").Value = 2.17689
$ws.Range("O13").Value = 0.2027419510312003
$ws.Range("P13").Value = 0.2027419510312003
$ws.Range("Q13").Value = 1.9556163878
$ws.Range("R13").Value = 17.6005474902
$ws.Range("S13").Value = 0.0001124051530142215
$ws.Range("T13").Value = 0.0001124051530142215

# Row 14: sCs -> ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col1a1"
$ws.Range("C14").Value = "Itga2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 28.08283533333333
$ws.Range("H14").Value = 84.24850599999999
$ws.Range("I14").Value = 0.005777169566755752
$ws.Range("J14").Value = 0.005777169566755752
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.740822
$ws.Range("N14").Value = 5.222466
$ws.Range("O14").Value = 0.4863878955914668
$ws.Range("P14").Value = 0.4863878955914669
$ws.Range("Q14").Value = 48.88721757064399
$ws.Range("R14").Value = 439.9849581357959
$ws.Range("S14").Value = 0.002809945348049396
$ws.Range("T14").Value = 0.002809945348049397

# Row 15: sCs -> FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col1a1"
$ws.Range("C15").Value = "Itga2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 28.08283533333333
$ws.Range("H15").Value = 84.24850599999999
$ws.Range("I15").Value = 0.005777169566755752
$ws.Range("J15").Value = 0.005777169566755752
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.077748
$ws.Range("N15").Value = 3.233244
$ws.Range("O15").Value = 0.3011241710513264
$ws.Range("P15").Value = 0.3011241710513265
$ws.Range("Q15").Value = 30.26621961482933
$ws.Range("R15").Value = 272.395976533464
$ws.Range("S15").Value = 0.001739645396812277
$ws.Range("T15").Value = 0.001739645396812277

# Row 16: sCs -> M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col1a1"
$ws.Range("C16").Value = "Itga2"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 28.08283533333333
$ws.Range("H16").Value = 84.24850599999999
$ws.Range("I16").Value = 0.005777169566755752
$ws.Range("J16").Value = 0.005777169566755752
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03488166666666666
$ws.Range("N16").Value = 0.104645
$ws.Range("O16").Value = 0.009745982326006345
$ws.Range("P16").Value = 0.009745982326006345
$ws.Range("Q16").Value = 0.9795761011522222
$ws.Range("R16").Value = 8.81618491037
$ws.Range("S16").Value = 0.0000563041924919433
$ws.Range("T16").Value = 0.0000563041924919433

# Row 17: sCs -> sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col1a1"
$ws.Range("C17").Value = "Itga2"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 28.08283533333333
$ws.Range("H17").Value = 84.24850599999999
$ws.Range("I17").Value = 0.005777169566755752
$ws.Range("J17").Value = 0.005777169566755752
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.7256300000000001
$ws.Range("N17").Value = 2.17689
$ws.Range("O17").Value = 0.2027419510312003
$ws.Range("P17").Value = 0.2027419510312003
$ws.Range("Q17").Value = 20.37774780292667
$ws.Range("R17").Value = 183.39973022634
$ws.Range("S17").Value = 0.001171274629402135
$ws.Range("T17").Value = 0.001171274629402135
